$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Row 6: the "_id" field becomes a plain "string" question (was previously a
#     special "db-object" appearance field with a repeat_count of 0). ---
# Match the formatting already used by the sibling cells in the row (B6:D6)
# instead of the special blue-highlighted font that only A6 had.
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(6, 1).Value = "string"

# appearance ("db-object") and repeat_count ("0") no longer apply - remove them
$ws.Cells.Item(6, 6).Clear()
$ws.Cells.Item(6, 12).Clear()

# --- Row 13 ("is_mch_instance" calculate field) is no longer needed. ---
$ws.Rows.Item(13).Delete()

# Selection moves to C6 to match the author's saved cursor position.
$ws.Range("C6").Select()
